$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Lembar1")

# Append two new normalisation pairs to the bottom of the table
$ws.Range("A119").Value = "cm"
$ws.Range("B119").Value = "Cuma"
$ws.Range("A120").Value = "g"
$ws.Range("B120").Value = "tidak"

# Match the selection state left by the edit
$ws.Range("C120").Select()
$excel.ActiveWindow.ScrollRow = 104
